$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 700
$ws.Range("B3").Value = 700
$ws.Range("B4").Value = 120
$ws.Range("B5").Value = 120
$ws.Range("B6").Value = 300
